$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: preserve the existing row 29 and row 30 data by copying them down
# to the new rows 31 and 32 before overwriting rows 29/30 with new values.
for ($col = 1; $col -le 20; $col++) {
    $src = $ws.Cells.Item(29, $col)
    $dst = $ws.Cells.Item(31, $col)
    $dst.Value = $src.Value2

    $src2 = $ws.Cells.Item(30, $col)
    $dst2 = $ws.Cells.Item(32, $col)
    $dst2.Value = $src2.Value2
}

# column D (date) keeps the same date number format as the source rows
$ws.Range("D31").NumberFormat = $ws.Range("D29").NumberFormat
$ws.Range("D32").NumberFormat = $ws.Range("D30").NumberFormat

# Step 2: update row 29 with the new (Primera) reading for Packham's Triumph
$ws.Range("D29").Value = 44747
$ws.Range("L29").Value = "Primera"
$ws.Range("N29").Value = 19000
$ws.Range("O29").Value = 20000
$ws.Range("P29").Value = 19500
$ws.Range("S29").Value = 1083

# Step 3: row 30 becomes a new record (Winter Nelis, Primera)
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44747
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104005
$ws.Range("J30").Value = "Pera"
$ws.Range("K30").Value = "Winter Nelis"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 19000
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 19500
$ws.Range("Q30").Value = "$/bandeja 18 kilos granel"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 1083
$ws.Range("T30").Value = 18
